# Weekly update: insert a new "Ajo" (garlic) price record for
# Macroferia Regional de Talca right before the existing row 419,
# shifting all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 419 (pushes old rows 419:522 to 420:523,
# extends the used range / dimension to A1:R523, and copies row formatting
# from the row above, matching Excel's native Insert behavior).
$ws.Rows.Item(419).Insert()

# Populate the newly inserted row 419 with the new weekly data point.
# (single-quoted literals so the "$/..." unit string is never interpolated)
$ws.Cells.Item(419, 1).Value = 5
$ws.Cells.Item(419, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(419, 3).Value = 'Maule'
$ws.Cells.Item(419, 4).Value = 45135
$ws.Cells.Item(419, 5).Value = 7
$ws.Cells.Item(419, 6).Value = 100112003
$ws.Cells.Item(419, 7).Value = 'Ajo'
$ws.Cells.Item(419, 8).Value = 'Chino'
$ws.Cells.Item(419, 9).Value = 'Primera'
$ws.Cells.Item(419, 10).Value = 200
$ws.Cells.Item(419, 11).Value = 20000
$ws.Cells.Item(419, 12).Value = 20000
$ws.Cells.Item(419, 13).Value = 20000
$ws.Cells.Item(419, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(419, 15).Value = 'China'
$ws.Cells.Item(419, 16).Value = 2000
$ws.Cells.Item(419, 17).Value = 10
$ws.Cells.Item(419, 18).Value = 'Hortaliza'
